# Generate Report for Handback
# Replace the two UUID-named source files (and their derived hash / Xliff
# file names / timestamps) that this handback run produced, across the
# Overview, zh-cn and de-de worksheets, including the cell values and the
# hyperlink "display" text that mirrors them.

$wb = $excel.ActiveWorkbook

$oldId1 = "362ff2f9-3b05-40e3-8c11-97be5941e3cb"
$newId1 = "866d8a5e-9b8d-490c-936d-5ac644844b30"
$oldId2 = "ac767584-6af0-470d-a3c7-014cc8455e3f"
$newId2 = "ffff87c334c6-3b43-4bbe-8282-6b5c7ee31d34"

$newHash = "3df4020b51e0a50d53e35174a1b02d4d3f489565"

$newXlfZh = "$newId1.$newHash.zh-cn.xlf"
$newXlfDe = "$newId1.$newHash.de-de.xlf"

$newOverviewDate = "2016-09-05 19:10:12"
$newZhHoDate = "2016-09-05 19:09:59"
$newZhHbDate = "2016-09-05 19:10:32"
$newDeHbDate = "2016-09-05 19:10:40"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId1.md"
$wsOverview.Range("B2").Value = "e2e\$newId1.md"
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Range("A3").Value = "$newId2.md"
$wsOverview.Range("B3").Value = "e2e\$newId2.md"
$wsOverview.Range("G3").Value = $newOverviewDate

# Hyperlinks only expose an "add new" primitive reliably, so rebuild the
# worksheet's hyperlinks from scratch, keeping the same target addresses
# and only updating the displayed text.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6af5b141368fd80643c2c0a33704092effd2456/e2e/$oldId1.md", "", "", "e2e\$newId1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6af5b141368fd80643c2c0a33704092effd2456/e2e/$oldId2.md", "", "", "e2e\$newId2.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newId1.md"
$wsZh.Range("I2").Value = "$newId1.md"
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("J2").Value = $newXlfZh
$wsZh.Range("H2").Value = $newZhHoDate
$wsZh.Range("K2").Value = $newZhHbDate

$wsZh.Range("A3").Value = "$newId2.md"
$wsZh.Range("I3").Value = "$newId2.md"
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("J3").Value = $newXlfZh
$wsZh.Range("H3").Value = $newZhHoDate
$wsZh.Range("K3").Value = $newZhHbDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6af5b141368fd80643c2c0a33704092effd2456/e2e/$oldId1.md", "", "", "$newId1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ae926d89f88d8eac50953191b051b035d1d50169/e2e/$oldId1.md", "", "", "$newId1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6af5b141368fd80643c2c0a33704092effd2456/e2e/$oldId2.md", "", "", "$newId2.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ae926d89f88d8eac50953191b051b035d1d50169/e2e/$oldId2.md", "", "", "$newId2.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newId1.md"
$wsDe.Range("I2").Value = "$newId1.md"
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("J2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newOverviewDate
$wsDe.Range("K2").Value = $newDeHbDate

$wsDe.Range("A3").Value = "$newId2.md"
$wsDe.Range("I3").Value = "$newId2.md"
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("J3").Value = $newXlfDe
$wsDe.Range("H3").Value = $newOverviewDate
$wsDe.Range("K3").Value = $newDeHbDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6af5b141368fd80643c2c0a33704092effd2456/e2e/$oldId1.md", "", "", "$newId1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8242edcbecb6b5f2d46baa7ac3863eb4441018bd/e2e/$oldId1.md", "", "", "$newId1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6af5b141368fd80643c2c0a33704092effd2456/e2e/$oldId2.md", "", "", "$newId2.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8242edcbecb6b5f2d46baa7ac3863eb4441018bd/e2e/$oldId2.md", "", "", "$newId2.md")
